$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 3 (positionId 2 / id 2) ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "P-Res"
$ws.Range("D3").Value = "P1 -P2, P2 -P3"
$ws.Range("E3").Value = "With people Responsible(EWB)`nWith people Responsible(IDN & EWB)"

# --- New row 4 (positionId 3 / id 3) ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "SPjM-SW`nSPjM-SW (without ECU-PjMresponsible in Department)`nECU-PjM"
$ws.Range("D4").Value = "P3-P4, P1-P2`nP2-P4"
$ws.Range("E4").Value = "Category C, D, E No people Responsible`nCategory C, D, E With people Responsible`nCategory C, D, E With people Responsible`nCategory B With people Responsible"

# --- New row 5 (positionId 4 / id 4) ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "SPjM-SW"

# Wrap text on the cells that hold multi-line content (matches style index 2 in the
# original workbook, which is the "wrapText" cell style already used by E2)
$ws.Range("E3").WrapText = $true
$ws.Range("C4:E4").WrapText = $true

# Row heights so the wrapped, multi-line text is fully visible
$ws.Rows.Item(3).RowHeight = 37.5
$ws.Rows.Item(4).RowHeight = 100

# Column widths widened to fit the new, longer role/project/responsibility text
$ws.Columns.Item(3).ColumnWidth = 21.666666666666664
$ws.Columns.Item(4).ColumnWidth = 25
$ws.Columns.Item(5).ColumnWidth = 25.833333333333336

# Selection moved by the author while reviewing the new rows
$ws.Range("D8").Select() | Out-Null
